$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Row 6: clear D6 (was "BMS APPS", now blank but keeps its style)
$ws.Range("D6").ClearContents()

# Row 11: D11 "SW" -> "APPS"; G11 becomes "T"
$ws.Range("D11").Value = "APPS"
$ws.Range("G11").Value = "T"

# Row 12: new "IVPDB" node - clone formatting from row 11, then fill in values
# (IVPDB string must be written before "BMS SW" so shared-string order matches)
$src = $ws.Range("C11:H11")
$dst = $ws.Range("C12:H12")
$src.Copy($dst)
$ws.Range("C12").Value = "IVPDB"
$ws.Range("D12").Value = "SW"
$ws.Range("F12").Value = " "
$ws.Range("G12").Value = "T"

# Row 7: D7 "SW" -> "BMS SW"; G7 "F" -> "T"
$ws.Range("D7").Value = "BMS SW"
$ws.Range("G7").Value = "T"

# Update active sheet/selection: "main" becomes the active tab with D7:D10 selected
$ws.Activate()
$ws.Range("D7:D10").Select()
